$wb = $excel.ActiveWorkbook

# Sheet "展览" updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1141
$ws1.Range("F3").Value = 634
$ws1.Range("F5").Value = 0
$ws1.Range("F7").Value = 9095
$ws1.Range("F11").Value = 622

# Sheet "全部类型" updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1141
$ws4.Range("F3").Value = 634
$ws4.Range("F7").Value = 4941
$ws4.Range("F10").Value = 9095
$ws4.Range("F16").Value = 622
